# SIMACITY POS+INV - APR 2020.xlsx
#
# Commit: "Misc fixes and enhancements ... Updated some files to function
# properly (removed rows with headers, etc)."
#
# The "20200430Iventory" sheet had 3 header/title rows (a big title
# "SiTime Inventory Format", a "Month" label, and a spacer row) sitting
# above the real table header. Those get removed so the sheet "functions
# properly" (e.g. so the real header row is row 1), and the inventory tab
# becomes the active/selected sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)            # "202004POS"
$ws2 = $wb.Worksheets.Item("20200430Iventory")

# --- Remove the three leading rows (title / month / blank spacer) -------
$ws2.Rows("1:3").Delete()

# --- The inventory sheet becomes the active / selected tab --------------
$ws2.Activate()

# Re-establish a header-row freeze (was frozen below the old row 4, now
# it should be frozen below the new row 1) and restore the selection.
$ws2.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $false
$excel.ActiveWindow.FreezePanes = $true
$ws2.Range("E3").Select()

# --- Refresh the AutoFilter over the new, smaller data range ------------
$ws2.AutoFilterMode = $false
$ws2.Range("A1:F10").AutoFilter()

# --- Update the _FilterDatabase defined name to match --------------------
$wb.Names.Item(1).RefersTo = "='20200430Iventory'!`$A`$1:`$F`$10"
